# Swap the content of column B and column C (rows 1-11) on the active sheet.
# This mirrors the source edit where the LastName/FirstName columns were
# swapped (B<->C) for the header row and all 10 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 11; $r++) {
    $bCell = $ws.Cells.Item($r, 2)   # column B
    $cCell = $ws.Cells.Item($r, 3)   # column C

    $bValue = $bCell.Value2
    $cValue = $cCell.Value2

    $bCell.Value2 = $cValue
    $cCell.Value2 = $bValue
}

# Update the sheet view: move the active selection from L5 to F8 (this also
# clears the stale topLeftCell="B1" scroll override from the old selection).
$ws.Range("F8").Select()
